# Workbook edit: merge relay rows RL18 + RL19-21 into a single "RL18-21" row
# on the Standalone ATU sheet, clear the reference for the LED row that used
# to be just below them, delete the now-empty spacer row, and fix a wording
# typo ("26SWG" -> "26 SWG") on the Embedded ATU sheet.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "Embedded ATU"
$ws2 = $wb.Worksheets.Item(2)   # "Standalone ATU"

# --- Standalone ATU: merge the RL18 and RL19-21 relay rows (13 & 14) ---
# Row 14 already has the correct final part/stock info (Omron relay), so
# keep that row's C/D/E values, just combine the quantities and reference.
$ws2.Range("A13").Value = 4
$ws2.Range("B13").Value = "RL18-21"
$ws2.Range("C13").Value = $ws2.Range("C14").Value2
$ws2.Range("D13").Value = $ws2.Range("D14").Value2
$ws2.Range("E13").Value = $ws2.Range("E14").Value2
$ws2.Rows.Item(13).RowHeight = 30

# --- Embedded ATU: fix wording of the VSWR bridge secondary winding note ---
$ws1.Range("C65").Value = "secondary: 10T 26 SWG wound around the sides of the core"

# Clear the reference text for the LED row (currently row 15); its other
# columns (qty/type/part) stay as-is.
$ws2.Range("B15").Value = ""

# Delete the old row 14 (the now-duplicated RL19-21 row); this shifts row
# 15 (LED row, with its reference already cleared) up to become row 14,
# and shifts everything below up by one as well.
$ws2.Rows.Item(14).Delete()

# --- Restore view/selection state seen in the saved workbook ---
$ws2.Activate()
$ws2.Rows.Item(13).Select() | Out-Null

$ws1.Activate()
$ws1.Range("F67").Select() | Out-Null
